# Generate Report for Handoff
# Updates the localization-status report:
#  - Refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#    for the 32ad6138-66d4-43bb-8a8e-d57d0c0a2d91 file across Overview/zh-cn/de-de sheets
#  - Sets the "Priority" column to "ht" (handoff type) for the newly-handed-off rows

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 13)

# -- Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-04 20:25:46"
}

# -- zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-04 20:25:41"
}

# -- de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-04 20:25:46"
}
